$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 28 — this shifts the existing rows 28-100
# down to 29-101 and grows the used range to A1:R101.
$ws.Rows.Item(28).Insert()

# Populate the newly inserted row 28 with the new weekly price record
# (same market/category/variety/quality/unit/origin as the row that used
# to be at 28, but a new date and its own volume/price figures).
$ws.Range("A28").Value = 7
$ws.Range("B28").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C28").Value = "Ñuble"
$ws.Range("D28").Value = 44414
$ws.Range("E28").Value = 16
$ws.Range("F28").Value = 100112003
$ws.Range("G28").Value = "Ajo"
$ws.Range("H28").Value = "Chino"
$ws.Range("I28").Value = "Primera"
$ws.Range("J28").Value = 120
$ws.Range("K28").Value = 14000
$ws.Range("L28").Value = 15000
$ws.Range("M28").Value = 14500
$ws.Range("N28").Value = "$/caja 10 kilos"
$ws.Range("O28").Value = "China"
$ws.Range("P28").Value = 1450
$ws.Range("Q28").Value = 10
$ws.Range("R28").Value = "Hortaliza"
